# Apply "post class 20 and 21" edit: Grades sheet gets a "Brief Submission"
# status column (P) and an assigned-citation column (Q) per student, a
# sheet-scoped "citation" defined name pointing at the example cell, a
# portrait page setup, and the Lab-Roster sheet gets a couple of cosmetic
# tweaks (selection + autosized columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grades")

# Header for the new "Brief Submission" status column.
$ws.Range("P1").Value = "Brief Submission"

$ws.Range("P2").Value = 'x'
$ws.Range("Q2").Value = 'Public Perception of Agricultural Pollution and Gulf of Mexico Hypoxia'
$ws.Range("P3").Value = 'x'
$ws.Range("Q3").Value = 'Implications of piscine predator control on the federally listed fountain darter. '
$ws.Range("P4").Value = 'x'
$ws.Range("Q4").Value = 'Relatedness and body size influence territorial behaviour in Salmo salar juveniles in the wild.'
$ws.Range("P5").Value = 'xx'
$ws.Range("P6").Value = 'x'
$ws.Range("Q6").Value = 'Growth response of largemouth bass (Micropterus salmoides) to catch-and-release angling: a 27-year mark–recapture study'
$ws.Range("P7").Value = 'x'
$ws.Range("Q7").Value = 'Accounting for variable recruitment and fishing mortality in 1 length-based stock assessments for data-limited fisheries'
$ws.Range("P8").Value = 'x'
$ws.Range("Q8").Value = 'Red Snapper Distribution on Natural Habitats and Artificial Structures in the Northern Gulf of Mexico'
$ws.Range("P9").Value = 'Late 2 hours'
$ws.Range("Q9").Value = 'Towards a balanced presentation and objective interpretation of acoustic and trawl survey data, with specific reference to the eastern Scotian Shelf'
$ws.Range("P10").Value = 'x'
$ws.Range("Q10").Value = 'Contrasting patterns of productivity and survival rates for stream-type chinook salmon (Oncorhynchus tshawytscha) populations of the Snake and Columbia rivers'
$ws.Range("P11").Value = 'x'
$ws.Range("Q11").Value = 'Management issues in the Lake Victoria watershed'
$ws.Range("P12").Value = 'x'
$ws.Range("Q12").Value = 'Using reverse-time egg transport analysis for predicting Asian carp spawning grounds in the Illinois River'
$ws.Range("P13").Value = 'x'
$ws.Range("Q13").Value = 'Comparing commercial and recreational harvest characteristics of paddlefish Polyodon spathula (Walbaum, 1792) in the Middle Mississippi River'
$ws.Range("P14").Value = 'x'
$ws.Range("Q14").Value = 'Temperature and hydrologic alteration predict the spread of invasive Largemouth Bass (Micropterus salmoides)'
$ws.Range("P15").Value = 'x'
$ws.Range("Q15").Value = 'Effects of hot dry summers on the loss of Atlantic salmon, Salmo salar, from estuaries in South West England'
$ws.Range("P16").Value = 'x'
$ws.Range("Q16").Value = 'GIS visualisation and analysis of mobile hydroacoustic fisheries data: a practical example'
$ws.Range("P18").Value = 'x'
$ws.Range("Q18").Value = 'Effectively managing angler satisfaction in recreational fisheries requires understanding the fish species and the anglers'
$ws.Range("P19").Value = 0
$ws.Range("P20").Value = 'x'
$ws.Range("Q20").Value = 'Tracking bowfin with acoustic telemetry: Insight into the ecology of a living fossil'
$ws.Range("P21").Value = 'x'
$ws.Range("Q21").Value = 'Assessing a social norms approach for improving recreational fisheries compliance'
$ws.Range("P22").Value = 'x'
$ws.Range("Q22").Value = 'Predictive Evaluation of Size Restrictions as Management Strategies for Tennessee Reservoir Crappie Fisheries'
$ws.Range("P23").Value = 'x'
$ws.Range("Q23").Value = 'The effects of fisheries management practises on freshwater ecosystems'
$ws.Range("P24").Value = 'x'
$ws.Range("Q24").Value = 'Effects of Multiple Low-Head Dams on Fish, Macroinvertebrates, Habitat, and Water Quality in the Fox River, Illinois'
$ws.Range("P25").Value = 'x'
$ws.Range("Q25").Value = 'Influence of behavior and mating success on brood-specific contribution to fish recruitment in ponds'

# Sheet-scoped defined name "citation" on Grades pointing at an example
# assigned-citation cell (Q3).
$ws.Names.Add("citation", "=Grades!`$Q`$3")

# Set this sheet's print orientation to portrait (matches target pageSetup).
$ws.PageSetup.Orientation = 1

# Move the active selection to where the instructor left off (row 22).
$ws.Range("P22").Select() | Out-Null

# Lab-Roster: minor view/format touch-ups that came along with the edit
# (columns E/F sized to fit their content).
$lr = $wb.Worksheets.Item("Lab-Roster")
$lr.Columns("E").ColumnWidth = 8.3
$lr.Columns("F").ColumnWidth = 7.6
$lr.Range("H5").Select() | Out-Null

# Leave Grades as the active/selected sheet (matches tabSelected in target).
$ws.Activate()
